$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 324.2
$ws.Cells.Item(2, 9).Value = 155.25
$ws.Cells.Item(2, 10).Value = 1000
$ws.Cells.Item(2, 11).Value = 155.25
$ws.Cells.Item(2, 12).Value = 1000
$ws.Cells.Item(2, 13).Value = -42.25
$ws.Cells.Item(2, 14).Value = -1226
$ws.Cells.Item(6, 8).Value = 1380.625
$ws.Cells.Item(6, 9).Value = 1549.2858
$ws.Cells.Item(6, 10).Value = 200
$ws.Cells.Item(6, 11).Value = 4647.857400000001
$ws.Cells.Item(6, 12).Value = 600
$ws.Cells.Item(6, 13).Value = -4535.857400000001
$ws.Cells.Item(6, 14).Value = -824
$ws.Cells.Item(21, 8).Value = 2283.3333
$ws.Cells.Item(21, 9).Value = 850
$ws.Cells.Item(21, 10).Value = 3000
$ws.Cells.Item(21, 11).Value = 850
$ws.Cells.Item(21, 12).Value = 3000
$ws.Cells.Item(21, 13).Value = -382
$ws.Cells.Item(21, 14).Value = -3936
$ws.Cells.Item(23, 8).Value = 2283.3333
$ws.Cells.Item(23, 9).Value = 850
$ws.Cells.Item(23, 10).Value = 3000
$ws.Cells.Item(23, 11).Value = 850
$ws.Cells.Item(23, 12).Value = 3000
$ws.Cells.Item(23, 13).Value = -616
$ws.Cells.Item(23, 14).Value = -3468
$ws.Cells.Item(137, 8).Value = 4423.8
$ws.Cells.Item(137, 9).Value = 3032.25
$ws.Cells.Item(137, 10).Value = 9990
$ws.Cells.Item(137, 11).Value = 9096.75
$ws.Cells.Item(137, 12).Value = 29970
$ws.Cells.Item(137, 13).Value = -6546.75
$ws.Cells.Item(137, 14).Value = -35070
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(36, 8).Value = 8666.666999999999
$ws.Cells.Item(36, 9).Value = 8666.666999999999
$ws.Cells.Item(36, 10).Value = 0
$ws.Cells.Item(36, 11).Value = 8666.666999999999
$ws.Cells.Item(36, 12).Value = 0
$ws.Cells.Item(36, 13).Value = -8320.666999999999
$ws.Cells.Item(122, 8).Value = 3688.5
$ws.Cells.Item(122, 9).Value = 3940.7778
$ws.Cells.Item(122, 10).Value = 2931.6667
$ws.Cells.Item(122, 11).Value = 11822.3334
$ws.Cells.Item(122, 12).Value = 8795.000100000001
$ws.Cells.Item(122, 13).Value = -9372.3334
$ws.Cells.Item(122, 14).Value = -13695.0001
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(5, 8).Value = 231.22223
$ws.Cells.Item(5, 9).Value = 196.83333
$ws.Cells.Item(5, 10).Value = 300
$ws.Cells.Item(5, 11).Value = 196.83333
$ws.Cells.Item(5, 12).Value = 300
$ws.Cells.Item(5, 13).Value = -83.83332999999999
$ws.Cells.Item(5, 14).Value = -526
$ws.Cells.Item(20, 8).Value = 1116
$ws.Cells.Item(20, 9).Value = 1150
$ws.Cells.Item(20, 10).Value = 1093.3334
$ws.Cells.Item(20, 11).Value = 1150
$ws.Cells.Item(20, 12).Value = 1093.3334
$ws.Cells.Item(20, 13).Value = -903
$ws.Cells.Item(20, 14).Value = -1587.3334
$ws.Cells.Item(26, 8).Value = 21030
$ws.Cells.Item(26, 9).Value = 21030
$ws.Cells.Item(26, 10).Value = 0
$ws.Cells.Item(26, 11).Value = 21030
$ws.Cells.Item(26, 12).Value = 0
$ws.Cells.Item(26, 13).Value = -20738
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 1500
$ws.Cells.Item(31, 9).Value = 1500
$ws.Cells.Item(31, 10).Value = 0
$ws.Cells.Item(31, 11).Value = 1500
$ws.Cells.Item(31, 12).Value = 0
$ws.Cells.Item(31, 13).Value = -1205
$ws.Cells.Item(34, 8).Value = 1500
$ws.Cells.Item(34, 9).Value = 1500
$ws.Cells.Item(34, 10).Value = 0
$ws.Cells.Item(34, 11).Value = 1500
$ws.Cells.Item(34, 12).Value = 0
$ws.Cells.Item(34, 13).Value = -1298
$ws.Cells.Item(58, 8).Value = 342.25
$ws.Cells.Item(58, 9).Value = 289.66666
$ws.Cells.Item(58, 10).Value = 500
$ws.Cells.Item(58, 11).Value = 289.66666
$ws.Cells.Item(58, 12).Value = 500
$ws.Cells.Item(58, 13).Value = -86.66665999999998
$ws.Cells.Item(58, 14).Value = -906
$ws.Cells.Item(136, 8).Value = 342.25
$ws.Cells.Item(136, 9).Value = 289.66666
$ws.Cells.Item(136, 10).Value = 500
$ws.Cells.Item(136, 11).Value = 868.9999799999999
$ws.Cells.Item(136, 12).Value = 1500
$ws.Cells.Item(136, 13).Value = 1681.00002
$ws.Cells.Item(136, 14).Value = -6600
$ws.Cells.Item(141, 8).Value = 122223
$ws.Cells.Item(141, 9).Value = 0
$ws.Cells.Item(141, 10).Value = 122223
$ws.Cells.Item(141, 11).Value = 0
$ws.Cells.Item(141, 12).Value = 122223
$ws.Cells.Item(141, 14).Value = -132583
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(23, 8).Value = 90.28570999999999
$ws.Cells.Item(23, 9).Value = 0
$ws.Cells.Item(23, 10).Value = 90.28570999999999
$ws.Cells.Item(23, 11).Value = 0
$ws.Cells.Item(23, 12).Value = 270.85713
$ws.Cells.Item(23, 14).Value = -740.85713
$ws.Cells.Item(34, 8).Value = 5698.4
$ws.Cells.Item(34, 9).Value = 0
$ws.Cells.Item(34, 10).Value = 5698.4
$ws.Cells.Item(34, 11).Value = 0
$ws.Cells.Item(34, 12).Value = 17095.2
$ws.Cells.Item(34, 14).Value = -17263.2
$ws.Cells.Item(55, 8).Value = 2373.1035
$ws.Cells.Item(55, 9).Value = 1330
$ws.Cells.Item(55, 10).Value = 2922.1052
$ws.Cells.Item(55, 11).Value = 3990
$ws.Cells.Item(55, 12).Value = 8766.3156
$ws.Cells.Item(55, 13).Value = -3813
$ws.Cells.Item(55, 14).Value = -9120.3156
$ws.Cells.Item(98, 8).Value = 22.5
$ws.Cells.Item(98, 9).Value = 20
$ws.Cells.Item(98, 10).Value = 25
$ws.Cells.Item(98, 11).Value = 60
$ws.Cells.Item(98, 12).Value = 75
$ws.Cells.Item(98, 13).Value = 1438
$ws.Cells.Item(98, 14).Value = -3071
$ws.Cells.Item(139, 8).Value = 5104.25
$ws.Cells.Item(139, 9).Value = 209.5
$ws.Cells.Item(139, 10).Value = 9999
$ws.Cells.Item(139, 11).Value = 628.5
$ws.Cells.Item(139, 12).Value = 29997
$ws.Cells.Item(139, 13).Value = 4511.5
$ws.Cells.Item(139, 14).Value = -40277
$ws.Cells.Item(140, 8).Value = 1150.4
$ws.Cells.Item(140, 9).Value = 1150.4
$ws.Cells.Item(140, 10).Value = 0
$ws.Cells.Item(140, 11).Value = 3451.2
$ws.Cells.Item(140, 12).Value = 0
$ws.Cells.Item(140, 13).Value = 1728.8
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 166672670
$ws.Cells.Item(70, 9).Value = 9000
$ws.Cells.Item(70, 10).Value = 500000000
$ws.Cells.Item(70, 11).Value = 9000
$ws.Cells.Item(70, 12).Value = 500000000
$ws.Cells.Item(70, 13).Value = -8730
$ws.Cells.Item(70, 14).Value = -500000540
$ws.Cells.Item(73, 8).Value = 166672670
$ws.Cells.Item(73, 9).Value = 9000
$ws.Cells.Item(73, 10).Value = 500000000
$ws.Cells.Item(73, 11).Value = 9000
$ws.Cells.Item(73, 12).Value = 500000000
$ws.Cells.Item(73, 13).Value = -8064
$ws.Cells.Item(73, 14).Value = -500001872
$ws.Cells.Item(75, 8).Value = 15000
$ws.Cells.Item(75, 9).Value = 0
$ws.Cells.Item(75, 10).Value = 15000
$ws.Cells.Item(75, 11).Value = 0
$ws.Cells.Item(75, 12).Value = 15000
$ws.Cells.Item(75, 14).Value = -16748
$ws.Cells.Item(78, 8).Value = 15000
$ws.Cells.Item(78, 9).Value = 0
$ws.Cells.Item(78, 10).Value = 15000
$ws.Cells.Item(78, 11).Value = 0
$ws.Cells.Item(78, 12).Value = 45000
$ws.Cells.Item(78, 14).Value = -53736
$ws.Cells.Item(102, 8).Value = 3773.25
$ws.Cells.Item(102, 9).Value = 3866.3333
$ws.Cells.Item(102, 10).Value = 3494
$ws.Cells.Item(102, 11).Value = 3866.3333
$ws.Cells.Item(102, 12).Value = 3494
$ws.Cells.Item(102, 13).Value = -2244.3333
$ws.Cells.Item(102, 14).Value = -6738
$ws.Cells.Item(114, 8).Value = 43500
$ws.Cells.Item(114, 9).Value = 0
$ws.Cells.Item(114, 10).Value = 43500
$ws.Cells.Item(114, 11).Value = 0
$ws.Cells.Item(114, 12).Value = 43500
$ws.Cells.Item(114, 14).Value = -52178
$ws.Cells.Item(122, 8).Value = 2457.4167
$ws.Cells.Item(122, 9).Value = 1443.2222
$ws.Cells.Item(122, 10).Value = 5500
$ws.Cells.Item(122, 11).Value = 4329.6666
$ws.Cells.Item(122, 12).Value = 16500
$ws.Cells.Item(122, 13).Value = -1879.6666
$ws.Cells.Item(122, 14).Value = -21400
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 730.6667
$ws.Cells.Item(22, 9).Value = 400
$ws.Cells.Item(22, 10).Value = 896
$ws.Cells.Item(22, 11).Value = 400
$ws.Cells.Item(22, 12).Value = 896
$ws.Cells.Item(22, 13).Value = -105
$ws.Cells.Item(22, 14).Value = -1486
$ws.Cells.Item(25, 8).Value = 1000
$ws.Cells.Item(25, 9).Value = 1000
$ws.Cells.Item(25, 10).Value = 1000
$ws.Cells.Item(25, 11).Value = 1000
$ws.Cells.Item(25, 12).Value = 1000
$ws.Cells.Item(25, 13).Value = -770
$ws.Cells.Item(25, 14).Value = -1460
$ws.Cells.Item(27, 8).Value = 730.6667
$ws.Cells.Item(27, 9).Value = 400
$ws.Cells.Item(27, 10).Value = 896
$ws.Cells.Item(27, 11).Value = 400
$ws.Cells.Item(27, 12).Value = 896
$ws.Cells.Item(27, 13).Value = -293
$ws.Cells.Item(27, 14).Value = -1110
$ws.Cells.Item(40, 8).Value = 721716.4399999999
$ws.Cells.Item(40, 9).Value = 4671
$ws.Cells.Item(40, 10).Value = 1259500.5
$ws.Cells.Item(40, 11).Value = 4671
$ws.Cells.Item(40, 12).Value = 1259500.5
$ws.Cells.Item(40, 13).Value = -4535
$ws.Cells.Item(40, 14).Value = -1259772.5
$ws.Cells.Item(122, 8).Value = 1900
$ws.Cells.Item(122, 9).Value = 1900
$ws.Cells.Item(122, 10).Value = 0
$ws.Cells.Item(122, 11).Value = 5700
$ws.Cells.Item(122, 12).Value = 0
$ws.Cells.Item(122, 13).Value = -3250
$ws.Cells.Item(136, 8).Value = 1063227.2
$ws.Cells.Item(136, 9).Value = 715784.3
$ws.Cells.Item(136, 10).Value = 2279277.5
$ws.Cells.Item(136, 11).Value = 2147352.9
$ws.Cells.Item(136, 12).Value = 6837832.5
$ws.Cells.Item(136, 13).Value = -2144802.9
$ws.Cells.Item(136, 14).Value = -6842932.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 3228.1428
$ws.Cells.Item(81, 9).Value = 3228.1428
$ws.Cells.Item(81, 10).Value = 0
$ws.Cells.Item(81, 11).Value = 6456.2856
$ws.Cells.Item(81, 12).Value = 0
$ws.Cells.Item(81, 13).Value = -5395.2856
$ws.Cells.Item(84, 8).Value = 3228.1428
$ws.Cells.Item(84, 9).Value = 3228.1428
$ws.Cells.Item(84, 10).Value = 0
$ws.Cells.Item(84, 11).Value = 32281.428
$ws.Cells.Item(84, 12).Value = 0
$ws.Cells.Item(84, 13).Value = -26977.428
$ws.Cells.Item(122, 8).Value = 2537.8
$ws.Cells.Item(122, 9).Value = 2537.8
$ws.Cells.Item(122, 10).Value = 0
$ws.Cells.Item(122, 11).Value = 7613.400000000001
$ws.Cells.Item(122, 12).Value = 0
$ws.Cells.Item(122, 13).Value = -5163.400000000001
$ws.Cells.Item(131, 8).Value = 68995
$ws.Cells.Item(131, 9).Value = 0
$ws.Cells.Item(131, 10).Value = 68995
$ws.Cells.Item(131, 11).Value = 0
$ws.Cells.Item(131, 12).Value = 68995
$ws.Cells.Item(131, 14).Value = -79075
$ws.Cells.Item(132, 8).Value = 913.3333
$ws.Cells.Item(132, 9).Value = 870
$ws.Cells.Item(132, 10).Value = 1000
$ws.Cells.Item(132, 11).Value = 2610
$ws.Cells.Item(132, 12).Value = 3000
$ws.Cells.Item(132, 13).Value = -80
$ws.Cells.Item(132, 14).Value = -8060
$ws.Cells.Item(136, 8).Value = 1272.9166
$ws.Cells.Item(136, 9).Value = 1203.2609
$ws.Cells.Item(136, 10).Value = 2875
$ws.Cells.Item(136, 11).Value = 3609.7827
$ws.Cells.Item(136, 12).Value = 8625
$ws.Cells.Item(136, 13).Value = -1059.7827
$ws.Cells.Item(136, 14).Value = -13725
